$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 3
$ws.Range("G8").Value = 1.57
$ws.Range("G9").Value = 2.15
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 3.2
$ws.Range("J9").Value = 2.88
$ws.Range("L9").Value = 3.75
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 11
$ws.Range("Q9").Value = 1.85
$ws.Range("R9").Value = 2
$ws.Range("U9").Value = 1.67
$ws.Range("V9").Value = 2.1
$ws.Range("X9").Value = 11
$ws.Range("Y9").Value = 9
$ws.Range("Z9").Value = 21
$ws.Range("AA9").Value = 17
$ws.Range("AC9").Value = 11
$ws.Range("AD9").Value = 6.5
$ws.Range("AH9").Value = 11
$ws.Range("AI9").Value = 17
$ws.Range("AK9").Value = 34
$ws.Range("AL9").Value = 23
$ws.Range("AM9").Value = 29
$ws.Range("AN9").Value = 4.33
$ws.Range("AO9").Value = 12
$ws.Range("AQ9").Value = 41
$ws.Range("AW9").Value = 5
$ws.Range("AX9").Value = 17
$ws.Range("AY9").Value = 23
$ws.Range("AZ9").Value = 51
$ws.Range("BA9").Value = 67
$ws.Range("G10").Value = 3.9
$ws.Range("I10").Value = 1.85
$ws.Range("J10").Value = 4.33
$ws.Range("L10").Value = 2.5
$ws.Range("N10").Value = 12
$ws.Range("S10").Value = 1.36
$ws.Range("T10").Value = 3
$ws.Range("U10").Value = 1.67
$ws.Range("V10").Value = 2.1
$ws.Range("X10").Value = 21
$ws.Range("Y10").Value = 13
$ws.Range("AA10").Value = 29
$ws.Range("AC10").Value = 12
$ws.Range("AE10").Value = 13
$ws.Range("AG10").Value = 151
$ws.Range("AH10").Value = 8.5
$ws.Range("AI10").Value = 9.5
$ws.Range("AK10").Value = 17
$ws.Range("AL10").Value = 15
$ws.Range("AO10").Value = 21
$ws.Range("AP10").Value = 26
$ws.Range("AT10").Value = 3
$ws.Range("AU10").Value = 7.5
$ws.Range("AW10").Value = 4
$ws.Range("AX10").Value = 10
$ws.Range("AZ10").Value = 34
$ws.Range("G14").Value = 1.91
$ws.Range("H14").Value = 3.5
$ws.Range("I14").Value = 3.8
$ws.Range("J14").Value = 2.6
$ws.Range("L14").Value = 4.5
$ws.Range("O14").Value = 1.3
$ws.Range("P14").Value = 3.4
$ws.Range("R14").Value = 1.8
$ws.Range("U14").Value = 1.91
$ws.Range("V14").Value = 1.91
$ws.Range("X14").Value = 9
$ws.Range("Y14").Value = 8.5
$ws.Range("Z14").Value = 17
$ws.Range("AH14").Value = 11
$ws.Range("AI14").Value = 19
$ws.Range("AJ14").Value = 13
$ws.Range("AL14").Value = 34
$ws.Range("AO14").Value = 11
$ws.Range("AP14").Value = 21
$ws.Range("AQ14").Value = 34
$ws.Range("AZ14").Value = 81
$ws.Range("BA14").Value = 101
$ws.Range("BB14").Value = 251
$ws.Range("Q17").Value = 2.15
$ws.Range("R17").Value = 1.67
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 5
$ws.Range("J18").Value = 7.5
$ws.Range("L18").Value = 1.8
$ws.Range("AE18").Value = 21
$ws.Range("AU18").Value = 9
$ws.Range("AX18").Value = 6
$ws.Range("H20").Value = 3.2
$ws.Range("M20").Value = 1.06
$ws.Range("N20").Value = 10
$ws.Range("Q20").Value = 1.93
$ws.Range("R20").Value = 1.93
$ws.Range("AC20").Value = 10
$ws.Range("AD20").Value = 6
$ws.Range("AJ20").Value = 12
$ws.Range("AM20").Value = 34
$ws.Range("G21").Value = 2.35
$ws.Range("J21").Value = 3.1
$ws.Range("L21").Value = 3.75
$ws.Range("N21").Value = 8.5
$ws.Range("S21").Value = 1.44
$ws.Range("T21").Value = 2.63
$ws.Range("AK21").Value = 34
$ws.Range("AO21").Value = 13
$ws.Range("AQ21").Value = 41
$ws.Range("AT21").Value = 2.63
$ws.Range("AW21").Value = 5
$ws.Range("K24").Value = 2.88
$ws.Range("M24").Value = 1.04
$ws.Range("N24").Value = 13
$ws.Range("Y24").Value = 12
$ws.Range("AF24").Value = 126
$ws.Range("AK24").Value = 301
$ws.Range("AL24").Value = 151
$ws.Range("G26").Value = 1.9
$ws.Range("H26").Value = 3.4
$ws.Range("I26").Value = 3.6
$ws.Range("J26").Value = 2.47
$ws.Range("K26").Value = 2.18
$ws.Range("L26").Value = 4.1
$ws.Range("M26").Value = 1.05
$ws.Range("N26").Value = 7.8
$ws.Range("O26").Value = 1.25
$ws.Range("P26").Value = 3.55
$ws.Range("Q26").Value = 1.75
$ws.Range("R26").Value = 2
$ws.Range("S26").Value = 1.36
$ws.Range("T26").Value = 2.87
$ws.Range("U26").Value = 1.65
$ws.Range("V26").Value = 2.12
$ws.Range("W26").Value = 8.5
$ws.Range("Z26").Value = 17.5
$ws.Range("AB26").Value = 22
$ws.Range("AC26").Value = 7.8
$ws.Range("AD26").Value = 6.7
$ws.Range("AE26").Value = 13
$ws.Range("AF26").Value = 50
$ws.Range("AG26").Value = 350
$ws.Range("AH26").Value = 11.75
$ws.Range("AI26").Value = 21
$ws.Range("AL26").Value = 32
$ws.Range("AM26").Value = 35
$ws.Range("AN26").Value = 3.95
$ws.Range("AP26").Value = 16.5
$ws.Range("AT26").Value = 2.87
$ws.Range("AU26").Value = 6.8
$ws.Range("AV26").Value = 55
$ws.Range("AX26").Value = 20
$ws.Range("AY26").Value = 25
$ws.Range("AZ26").Value = 110
$ws.Range("BA26").Value = 120
$ws.Range("BB26").Value = 300
